$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that used to sit at the end of
#    the "you currently use?" paragraph.
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2) Append four new paragraphs right after the "Confusion Matrix…"
#    paragraph (the last paragraph in the body, just before sectPr):
#      - an empty spacer paragraph
#      - a bold "Model Evaluation" heading
#      - a body paragraph (with the _GoBack bookmark re-added on it)
#      - a bold "Model Validation" heading
#    We build the insertion point as the position immediately before
#    the final paragraph mark of the document, so InsertXML appends
#    after the existing last paragraph instead of merging into it.
# ------------------------------------------------------------------
$endRange = $d.Content
$insertAt = $d.Range($endRange.End - 1, $endRange.End - 1)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPrPlain = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPrBold  = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$xml = ''
$xml += "<w:p $ns><w:pPr>$rPrPlain</w:pPr></w:p>"
$xml += "<w:p $ns><w:pPr>$rPrBold</w:pPr><w:r>$rPrBold<w:t>Model Evaluation</w:t></w:r></w:p>"
$xml += "<w:p $ns><w:pPr>$rPrPlain</w:pPr><w:r>$rPrPlain<w:t>As the dataset is perfectly balanced, we fixed the accuracy threshold to be 0.5.</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$xml += "<w:p $ns><w:pPr>$rPrBold</w:pPr><w:r>$rPrBold<w:t>Model Validation</w:t></w:r></w:p>"

$insertAt.InsertXML($xml)

Write-Output "done"
